# Fixing geothermal and ammonia issue, and time period of weather data and
# demand schedule.
#
# The "Plant interest/lifetime" and "Infrastructure interest/lifetime"
# columns (previously J:M) shift two columns to the right (to L:O) to make
# room for two new "Geothermal interest rate" / "Geothermal lifetime
# (years)" columns at J:K. The previously-blank "Hydro interest rate" /
# "Hydro lifetime (years)" values on row 2 (H2:I2) are also populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at J:K, shifting the old J:M ("Plant ...",
# "Infrastructure ...") columns right to L:O.
$ws.Range("J1:K2").Insert(-4161)

# Fill in the previously-empty Hydro values on row 2.
$ws.Range("H2").Value = 0.10881498793561939
$ws.Range("I2").Value = 20

# New Geothermal columns (J:K).
$ws.Range("J1").Value = "Geothermal interest rate"
$ws.Range("K1").Value = "Geothermal lifetime (years)"
$ws.Range("J2").Value = 0.10881498793561939
$ws.Range("K2").Value = 20

# Re-selection, matching the saved workbook view after editing.
$ws.Range("K5").Select() | Out-Null

# Resize columns to (approximately) match the author's final, content-fit
# column widths.
$ws.Columns.Item(1).ColumnWidth = 7.333333333333333
$ws.Columns.Item(2).ColumnWidth = 24
$ws.Columns.Item(3).ColumnWidth = 19.666666666666668
$ws.Columns.Item(4).ColumnWidth = 15.166666666666666
$ws.Columns.Item(5).ColumnWidth = 17.833333333333336
$ws.Columns.Item(6).ColumnWidth = 15.499999999999998
$ws.Columns.Item(7).ColumnWidth = 18
$ws.Columns.Item(8).ColumnWidth = 15.833333333333334
$ws.Columns.Item(9).ColumnWidth = 18.333333333333336
$ws.Columns.Item(10).ColumnWidth = 15.166666666666666
$ws.Columns.Item(11).ColumnWidth = 17.833333333333336
$ws.Columns.Item(12).ColumnWidth = 22.5
$ws.Columns.Item(13).ColumnWidth = 25
